$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without letting Excel
# auto-coerce numeric-looking strings (e.g. "230.44") into numbers, and
# without leaving the cells style pointing at a new "Text" format.
function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

# Refreshed "Price" (D) and "Volume(1h)" (E) columns from the latest crypto feed.
Set-TextValue $ws 'D2' '35.009.95'
Set-TextValue $ws 'E2' '  -0.22%  '
Set-TextValue $ws 'D3' '1.818.53'
Set-TextValue $ws 'E3' '  -0.25%  '
Set-TextValue $ws 'E4' '  -0.15%  '
Set-TextValue $ws 'D5' '230.44'
Set-TextValue $ws 'E5' '  -0.68%  '
Set-TextValue $ws 'E6' '  +0.56%  '
Set-TextValue $ws 'E7' '  -0.09%  '
Set-TextValue $ws 'D8' '40.05'
Set-TextValue $ws 'E8' '  -4.02%  '
Set-TextValue $ws 'D9' '0.324'
Set-TextValue $ws 'E9' '  +4.95%  '
Set-TextValue $ws 'D10' '0.0683'
Set-TextValue $ws 'E10' '  -0.24%  '
Set-TextValue $ws 'D11' '0.0992'
Set-TextValue $ws 'E11' '  -1.26%  '
Set-TextValue $ws 'D12' '2.081.58'
Set-TextValue $ws 'E12' '  -0.30%  '
Set-TextValue $ws 'D13' '11.32'
Set-TextValue $ws 'E13' '  +2.15%  '
Set-TextValue $ws 'D14' '0.667'
Set-TextValue $ws 'E14' '  +1.03%  '
Set-TextValue $ws 'D15' '1.816.46'
Set-TextValue $ws 'E15' '  -0.29%  '
Set-TextValue $ws 'D16' '4.62'
Set-TextValue $ws 'E16' '  -0.93%  '
Set-TextValue $ws 'D17' '34.989.12'
Set-TextValue $ws 'E17' '  -0.18%  '
Set-TextValue $ws 'E18' '  -0.17%  '
Set-TextValue $ws 'D19' '0.0₃0786'
Set-TextValue $ws 'E19' '  -0.22%  '
Set-TextValue $ws 'D20' '240.78'
Set-TextValue $ws 'E20' '  +0.82%  '
Set-TextValue $ws 'D21' '12.14'
Set-TextValue $ws 'E21' '  +3.23%  '
Set-TextValue $ws 'D22' '4.66'
Set-TextValue $ws 'E22' '  +1.57%  '
Set-TextValue $ws 'E23' '  -0.04%  '
Set-TextValue $ws 'D24' '2.26'
Set-TextValue $ws 'E24' '  +1.48%  '
Set-TextValue $ws 'D25' '173.46'
Set-TextValue $ws 'E25' '  +0.95%  '
Set-TextValue $ws 'D26' '7.84'
Set-TextValue $ws 'E26' '  +0.94%  '
Set-TextValue $ws 'D27' '0.124'
Set-TextValue $ws 'E27' '  +2.53%  '
Set-TextValue $ws 'D28' '17.35'
Set-TextValue $ws 'E28' '  -0.70%  '
Set-TextValue $ws 'E29' '  -4.98%  '
Set-TextValue $ws 'E30' '  -0.19%  '
Set-TextValue $ws 'E31' '  +2.45%  '
Set-TextValue $ws 'D32' '0.0550'
Set-TextValue $ws 'E32' '  -0.36%  '
Set-TextValue $ws 'D33' '3.95'
Set-TextValue $ws 'E33' '  -0.41%  '
Set-TextValue $ws 'D34' '1.25'
Set-TextValue $ws 'E34' '  +12.27%  '
Set-TextValue $ws 'D35' '1.83'
Set-TextValue $ws 'E35' '  +2.85%  '
Set-TextValue $ws 'D36' '0.691'
Set-TextValue $ws 'E36' '  +2.19%  '
Set-TextValue $ws 'D37' '92.45'
Set-TextValue $ws 'E37' '  -0.24%  '
Set-TextValue $ws 'E38' '  +5.58%  '
Set-TextValue $ws 'D39' '1.339.80'
Set-TextValue $ws 'E39' '  +1.79%  '
Set-TextValue $ws 'D40' '0.0194'
Set-TextValue $ws 'E40' '  +0.72%  '
Set-TextValue $ws 'D41' '0.981'
Set-TextValue $ws 'E41' '  -0.73%  '
Set-TextValue $ws 'D42' '14.63'
Set-TextValue $ws 'E42' '  -0.24%  '
Set-TextValue $ws 'E43' '  -2.82%  '
Set-TextValue $ws 'E44' '  -1.21%  '
Set-TextValue $ws 'D45' '2.75'
Set-TextValue $ws 'E45' '  -1.10%  '
Set-TextValue $ws 'E46' '  +2.04%  '
Set-TextValue $ws 'D47' '6.21'
Set-TextValue $ws 'E47' '  +0.62%  '
Set-TextValue $ws 'D48' '1.997.82'
Set-TextValue $ws 'E48' '  -0.15%  '
Set-TextValue $ws 'E49' '  -0.08%  '
Set-TextValue $ws 'E50' '  +3.44%  '
Set-TextValue $ws 'D51' '97.07'
Set-TextValue $ws 'E51' '  -3.08%  '
